$d = $word.ActiveDocument
$green = 5287936  # RGB(0, 176, 80) -> 00B050

# 1) "Make cases for why studies with school-aged students..." paragraph:
#    add green color to the paragraph mark and its run.
$p12 = $d.Paragraphs(12)
$p12.Range.Font.Color = $green

# 2) Merge the two runs of the "Grades, dropout, graduation, ..." paragraph
#    into a single run (removing the _GoBack bookmark that used to sit
#    between them).
$bm = $d.Bookmarks("_GoBack")
$splitPos = $bm.Start
$p17 = $d.Paragraphs(17)
$paraEnd = $p17.Range.End - 1
$r2 = $d.Range($splitPos, $paraEnd)
$r2text = $r2.Text
$bm.Delete()
$r2b = $d.Range($splitPos, $paraEnd)
$r2b.Delete()
$r1 = $d.Range($p17.Range.Start, $splitPos)
$r1.Text = $r1.Text + $r2text

# 3) "Sleep hygiene" paragraph: add green color to paragraph mark and run.
$p28 = $d.Paragraphs(28)
$p28.Range.Font.Color = $green

# 4) "ABCs of sleeping (make this consistent with the used measure)"
#    paragraph: add green color to paragraph mark and run.
$p29 = $d.Paragraphs(29)
$p29.Range.Font.Color = $green

# 5) Re-create the _GoBack bookmark spanning from the start of the
#    "Sleep hygiene" paragraph through the end of the "ABCs of sleeping..."
#    paragraph's text (before its paragraph mark).
$startPos = $p28.Range.Start
$endPos = $p29.Range.End - 1
$bmr = $d.Range($startPos, $endPos)
$bmr.Bookmarks.Add("_GoBack")
